$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.126.20'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.046.74'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.93%  '

$ws.Range("E6").Value = '  -1.49%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.79'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -4.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.382'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.40%  '

$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.892'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +9.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.346.15'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.71'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.048.92'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.39'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +11.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.153.61'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '74.45'
$ws.Range("D19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0894'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.42'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.71'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("E24").Value = '  +2.94%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.58'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.95'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.44%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.18'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.99%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.13'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.05%  '

$ws.Range("E29").Value = '  -1.36%  '

$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.16'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.78%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.93'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0620'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.48'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0883'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.25'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.28%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.34'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.18%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.27'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +14.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.10'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +8.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0988'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -15.83%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0223'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.57%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.54'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.30%  '

$ws.Range("E44").Value = '  -1.64%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '95.77'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.44'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.271.35'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.37%  '

$ws.Range("E48").Value = '  -2.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.78'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.231.62'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '44.31'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.05%  '
